# Add a new "Save" column (H) to the s_vals sheet, matching the header
# style used by the existing columns (B1:G1) and filling the data rows
# with 0, mirroring the other numeric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, copying formatting from the neighboring header (G1)
# so it keeps the same bold/bordered/centered style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data values for the "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
